$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: corrected gas price
# Leading apostrophe forces the value to stay literal text (matches the
# original inlineStr cell) instead of Excel auto-coercing it to a number;
# ClearFormats() strips the quote-prefix style Excel applies so the cell
# keeps using the default style, same as before.
$ws.Range("C2").Value = "'2.689"
$ws.Range("C2").ClearFormats()

# Row 3: product name + price were found this time
$ws.Range("B3").Value = "That's Smart! Fat Free Skim Milk"
$ws.Range("C3").Value = "'$2.72"
$ws.Range("C3").ClearFormats()

# Drop all the stale rows (6-20) from the previous run - only the first
# five rows (header + 4 data rows) survive.
$ws.Range("A6:C20").EntireRow.Delete()
